$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Uncomment RAD Extension Payments test data: change the Execute flag
# for the "Extension Payments" row (C4) from "DONOTRUN" to "Y".
$ws.Range("C4").Value = "Y"
